$wb = $excel.ActiveWorkbook

# Sheet "展览" updates (column F = 想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 170
$ws1.Range("F5").Value = 178
$ws1.Range("F6").Value = 2703
$ws1.Range("F7").Value = 187
$ws1.Range("F9").Value = 180
$ws1.Range("F10").Value = 1570
$ws1.Range("F11").Value = 547
$ws1.Range("F12").Value = 46
$ws1.Range("F24").Value = 1736
$ws1.Range("F27").Value = 70
$ws1.Range("F30").Value = 308

# Sheet "全部类型" updates (column F = 想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 170
$ws4.Range("F6").Value = 178
$ws4.Range("F7").Value = 2703
$ws4.Range("F8").Value = 187
$ws4.Range("F10").Value = 180
$ws4.Range("F11").Value = 1570
$ws4.Range("F12").Value = 547
$ws4.Range("F13").Value = 46
$ws4.Range("F25").Value = 1736
$ws4.Range("F28").Value = 70
$ws4.Range("F31").Value = 308
